# Update the "想去人数" (F column) counts on both the "展览" and
# "全部类型" worksheets, which carry identical data tables.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 2159
    6  = 12511
    7  = 12511
    8  = 57
    11 = 453
    14 = 13637
    15 = 13917
    20 = 1044
    23 = 441
    24 = 5018
    25 = 242
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
